# "a new quest of thief"
#
# Adds a new drop-table entry (row) for a new quest "盗贼II" ("Thief II")
# at the end of the "Drop" sheet's table (表2), mirroring the shape of the
# existing quest rows (e.g. row 44, "一个愿望").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a single table covering A3:I44 - grow it by one row so the
# table range, autofilter and dimension all extend to row 45 automatically.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()
$newRow = $newListRow.Range.Row

# Copy the formatting of the previous data row (44) into the freshly added
# row so the new cells pick up the same styles (B/C plain text, D/E wrapped
# item-list styles, I count style) instead of being left unformatted.
$prevRow = $newRow - 1
$ws.Range("B" + $prevRow).Copy()
$ws.Range("B" + $newRow).PasteSpecial(-4122)
$ws.Range("C" + $prevRow).Copy()
$ws.Range("C" + $newRow).PasteSpecial(-4122)
$ws.Range("D" + $prevRow).Copy()
$ws.Range("D" + $newRow).PasteSpecial(-4122)
$ws.Range("E" + $prevRow).Copy()
$ws.Range("E" + $newRow).PasteSpecial(-4122)
$ws.Range("I" + $prevRow).Copy()
$ws.Range("I" + $newRow).PasteSpecial(-4122)

# Fill in the new quest's data.
$ws.Range("A" + $newRow).Value = 23000505
$ws.Range("B" + $newRow).Value = "盗贼II"
$ws.Range("C" + $newRow).Value = "dlthief"
$ws.Range("D" + $newRow).Value = "suijihuanshouka;suijihuanshouka;sucaidai(an);xiaoxinghuoliyaoji"
$ws.Range("E" + $newRow).Value = "15;15;30;40"
$ws.Range("I" + $newRow).Value = 2

# Put the selection where Excel would naturally leave it after entering the
# last value of the new row.
$ws.Range("C" + $newRow).Select()
